$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.906.43'
$ws.Range('E2').Value = '  +2.26%  '
$ws.Range('D3').Value = '3.160.93'
$ws.Range('E3').Value = '  +4.17%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'579.26"
$ws.Range('E5').Value = '  +4.60%  '
$ws.Range('D6').Value = "'150.35"
$ws.Range('E6').Value = '  +7.12%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.158.42'
$ws.Range('E8').Value = '  +4.17%  '
$ws.Range('E9').Value = '  +2.12%  '
$ws.Range('E10').Value = '  +6.73%  '
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('D12').Value = "'0.501"
$ws.Range('E12').Value = '  +3.66%  '
$ws.Range('D13').Value = "'0.0000268"
$ws.Range('E13').Value = '  +17.87%  '
$ws.Range('D14').Value = "'37.52"
$ws.Range('E14').Value = '  +6.41%  '
$ws.Range('D15').Value = '3.676.28'
$ws.Range('E15').Value = '  +4.03%  '
$ws.Range('D16').Value = '64.950.72'
$ws.Range('E16').Value = '  +2.21%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.157.37'
$ws.Range('E17').Value = '  +4.03%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = "'7.15"
$ws.Range('E18').Value = '  +6.35%  '
$ws.Range('E19').Value = '  +1.65%  '
$ws.Range('D20').Value = "'510.07"
$ws.Range('E20').Value = '  +8.21%  '
$ws.Range('D21').Value = "'14.81"
$ws.Range('E21').Value = '  +5.95%  '
$ws.Range('E22').Value = '  +7.03%  '
$ws.Range('D23').Value = "'15.30"
$ws.Range('E23').Value = '  +6.01%  '
$ws.Range('E24').Value = '  +4.15%  '
$ws.Range('D25').Value = "'84.95"
$ws.Range('E25').Value = '  +3.16%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = "'9.01"
$ws.Range('E27').Value = '  +12.16%  '
$ws.Range('E28').Value = '  +5.28%  '
$ws.Range('E29').Value = '  +8.03%  '
$ws.Range('E30').Value = '  +6.84%  '
$ws.Range('E31').Value = '  +14.98%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').Value = '  +4.48%  '
$ws.Range('E34').Value = '  +11.68%  '
$ws.Range('D35').Value = "'6.57"
$ws.Range('E35').Value = '  +6.93%  '
$ws.Range('D36').Value = "'55.74"
$ws.Range('E36').Value = '  +2.03%  '
$ws.Range('E37').Value = '  +10.67%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = "'3.13"
$ws.Range('E38').Value = '  +14.02%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = "'471.24"
$ws.Range('E39').Value = '  +7.49%  '
$ws.Range('D40').Value = "'0.0421"
$ws.Range('E40').Value = '  +3.96%  '
$ws.Range('D41').Value = "'8.63"
$ws.Range('E41').Value = '  +4.84%  '
$ws.Range('D42').Value = '3.064.48'
$ws.Range('E42').Value = '  +2.56%  '
$ws.Range('E43').Value = '  +1.36%  '
$ws.Range('E44').Value = '  +6.07%  '
$ws.Range('E45').Value = '  +8.88%  '
$ws.Range('D46').Value = "'29.30"
$ws.Range('E46').Value = '  +6.32%  '
$ws.Range('D47').Value = '0.0₃0598'
$ws.Range('E47').Value = '  +17.98%  '
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('D50').Value = "'2.25"
$ws.Range('E50').Value = '  +9.20%  '
$ws.Range('D51').Value = "'119.68"
$ws.Range('E51').Value = '  +1.57%  '
